$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.250258139219767
$ws.Range("D2").Value = 0.005747792639855476
$ws.Range("E2").Value = 0.5268423054473814
$ws.Range("F2").Value = 0.9530285815336583
$ws.Range("G2").Value = 0.8715018868364552
$ws.Range("H2").Value = 0.7425514595447567
$ws.Range("L2").Value = 0.374049089480053
$ws.Range("N2").Value = 1.468496641925682
$ws.Range("B3").Value = 1.188230561460585
$ws.Range("D3").Value = 0.005663307123061401
$ws.Range("E3").Value = 0.5029502722519297
$ws.Range("F3").Value = 0.8946491321224102
$ws.Range("G3").Value = 0.8044451341010017
$ws.Range("H3").Value = 0.7175415177758282
$ws.Range("L3").Value = 0.3325800823766087
$ws.Range("N3").Value = 1.4692426749845
$ws.Range("B4").Value = 1.150682959045383
$ws.Range("D4").Value = 0.005621275071270304
$ws.Range("E4").Value = 0.4882034792257812
$ws.Range("F4").Value = 0.8595645133386114
$ws.Range("G4").Value = 0.763963207706098
$ws.Range("H4").Value = 0.7027713160572375
$ws.Range("L4").Value = 0.3071694250850783
$ws.Range("N4").Value = 1.47037184792687
$ws.Range("B5").Value = 1.1355174883291
$ws.Range("D5").Value = 0.005606577100008892
$ws.Range("E5").Value = 0.4821748288301819
$ws.Range("F5").Value = 0.8454561265490099
$ws.Range("G5").Value = 0.747637448999825
$ws.Range("H5").Value = 0.6968985148211573
$ws.Range("L5").Value = 0.2968272395221163
$ws.Range("N5").Value = 1.471001338671897
$ws.Range("B6").Value = 1.133007465481853
$ws.Range("D6").Value = 0.00560428200867058
$ws.Range("E6").Value = 0.4811726200858359
$ws.Range("F6").Value = 0.8431247784158131
$ws.Range("G6").Value = 0.7449368115371442
$ws.Range("H6").Value = 0.6959321350483663
$ws.Range("L6").Value = 0.2951107022479107
$ws.Range("N6").Value = 1.471116110577725
$ws.Range("B7").Value = 1.150477883242388
$ws.Range("D7").Value = 0.005621067067572483
$ws.Range("E7").Value = 0.4881222523195632
$ws.Range("F7").Value = 0.8593734810721969
$ws.Range("G7").Value = 0.7637423443791533
$ws.Range("H7").Value = 0.7026915231134012
$ws.Range("L7").Value = 0.3070298947698689
$ws.Range("N7").Value = 1.470379651012877
$ws.Range("B8").Value = 1.228759611986476
$ws.Range("D8").Value = 0.005716592311401314
$ws.Range("E8").Value = 0.5186203908175457
$ws.Range("F8").Value = 0.9327400290007546
$ws.Range("G8").Value = 0.848235519657436
$ws.Range("H8").Value = 0.733805684110564
$ws.Range("L8").Value = 0.3597398267595793
$ws.Range("N8").Value = 1.46861493217645
$ws.Range("B9").Value = 1.386531631542027
$ws.Range("D9").Value = 0.005984045202204413
$ws.Range("E9").Value = 0.5778128820736441
$ws.Range("F9").Value = 1.082762504600055
$ws.Range("G9").Value = 1.019548510353474
$ws.Range("H9").Value = 0.7995277220165349
$ws.Range("L9").Value = 0.4635216917293405
$ws.Range("N9").Value = 1.470455404179191
$ws.Range("B10").Value = 1.505054176528631
$ws.Range("D10").Value = 0.006232242572785651
$ws.Range("E10").Value = 0.6209270957724016
$ws.Range("F10").Value = 1.196906630558544
$ws.Range("G10").Value = 1.149043784464141
$ws.Range("H10").Value = 0.850769447951933
$ws.Range("L10").Value = 0.5400461550407272
$ws.Range("N10").Value = 1.475010192437068
$ws.Range("B11").Value = 1.559542429807664
$ws.Range("D11").Value = 0.00635700432626507
$ws.Range("E11").Value = 0.6404599221318534
$ws.Range("F11").Value = 1.249723867037801
$ws.Range("G11").Value = 1.208788306643726
$ws.Range("H11").Value = 0.8747415991251444
$ws.Range("L11").Value = 0.5749242202579978
$ws.Range("N11").Value = 1.477771919955529
$ws.Range("B12").Value = 1.580257902491269
$ws.Range("D12").Value = 0.006406004435667967
$ws.Range("E12").Value = 0.6478449279377685
$ws.Range("F12").Value = 1.26985562414535
$ws.Range("G12").Value = 1.231535725413437
$ws.Range("H12").Value = 0.8839158848437592
$ws.Range("L12").Value = 0.5881414874193069
$ws.Range("N12").Value = 1.478916396014441
$ws.Range("B13").Value = 1.575792813939302
$ws.Range("D13").Value = 0.006395372492534079
$ws.Range("E13").Value = 0.6462549551843324
$ws.Range("F13").Value = 1.265514020131093
$ws.Range("G13").Value = 1.226631119608726
$ws.Range("H13").Value = 0.8819357223336795
$ws.Range("L13").Value = 0.5852944791634229
$ws.Range("N13").Value = 1.478665532284282
$ws.Range("B14").Value = 1.561245063201625
$ws.Range("D14").Value = 0.006361000118417337
$ws.Range("E14").Value = 0.6410677253016104
$ws.Range("F14").Value = 1.251377477616984
$ws.Range("G14").Value = 1.210657257706799
$ws.Range("H14").Value = 0.8754944302671959
$ws.Range("L14").Value = 0.5760114177041942
$ws.Range("N14").Value = 1.477864102093861
$ws.Range("B15").Value = 1.552344811013199
$ws.Range("D15").Value = 0.006340176186363067
$ws.Range("E15").Value = 0.6378888757731005
$ws.Range("F15").Value = 1.242735579239223
$ws.Range("G15").Value = 1.200888979850248
$ws.Range("H15").Value = 0.8715615673102945
$ws.Range("L15").Value = 0.5703265412523422
$ws.Range("N15").Value = 1.477386037670925
$ws.Range("B16").Value = 1.501504723318533
$ws.Range("D16").Value = 0.006224331984398646
$ws.Range("E16").Value = 0.6196489566787307
$ws.Range("F16").Value = 1.193473087308803
$ws.Range("G16").Value = 1.145156442517873
$ws.Range("H16").Value = 0.8492162485185872
$ws.Range("L16").Value = 0.53776814278541
$ws.Range("N16").Value = 1.47484354209513
$ws.Range("B17").Value = 1.470462224042251
$ws.Range("D17").Value = 0.006156337701064274
$ws.Range("E17").Value = 0.608438752781737
$ws.Range("F17").Value = 1.16348262850903
$ws.Range("G17").Value = 1.111182828922381
$ws.Range("H17").Value = 0.8356786385495809
$ws.Range("L17").Value = 0.5178117786446137
$ws.Range("N17").Value = 1.473460091412775
$ws.Range("B18").Value = 1.452661230544209
$ws.Range("D18").Value = 0.006118342013934353
$ws.Range("E18").Value = 0.6019834137033229
$ws.Range("F18").Value = 1.146316754359532
$ws.Range("G18").Value = 1.091720645355707
$ws.Range("H18").Value = 0.8279543950656318
$ws.Range("L18").Value = 0.5063396771016926
$ws.Range("N18").Value = 1.472729317848518
$ws.Range("B19").Value = 1.44664337379038
$ws.Range("D19").Value = 0.006105666921747854
$ws.Range("E19").Value = 0.5997964585859563
$ws.Range("F19").Value = 1.140519003701968
$ws.Range("G19").Value = 1.085144477448097
$ws.Range("H19").Value = 0.8253497463661006
$ws.Range("L19").Value = 0.5024564945763075
$ws.Range("N19").Value = 1.472493061635191
$ws.Range("B20").Value = 1.473761181728833
$ws.Range("D20").Value = 0.006163460269448251
$ws.Range("E20").Value = 0.6096328781113698
$ws.Range("F20").Value = 1.166666467179169
$ws.Range("G20").Value = 1.114791226524432
$ws.Range("H20").Value = 0.8371132897556208
$ws.Range("L20").Value = 0.5199355178616827
$ws.Range("N20").Value = 1.473600643706405
$ws.Range("B21").Value = 1.565515868070122
$ws.Range("D21").Value = 0.006371048065268781
$ws.Range("E21").Value = 0.64259165776933
$ws.Range("F21").Value = 1.255526145610702
$ws.Range("G21").Value = 1.215345791828781
$ws.Range("H21").Value = 0.8773837626106911
$ws.Range("L21").Value = 0.5787378136850521
$ws.Range("N21").Value = 1.478096827723164
$ws.Range("B22").Value = 1.625960435985178
$ws.Range("D22").Value = 0.006516975326011476
$ws.Range("E22").Value = 0.6640641507683114
$ws.Range("F22").Value = 1.314366008253927
$ws.Range("G22").Value = 1.281785162844585
$ws.Range("H22").Value = 0.9042662663682961
$ws.Range("L22").Value = 0.6172251146049064
$ws.Range("N22").Value = 1.481610237635692
$ws.Range("B23").Value = 1.593656417867976
$ws.Range("D23").Value = 0.00643813579005581
$ws.Range("E23").Value = 0.6526101415744563
$ws.Range("F23").Value = 1.282891182970758
$ws.Range("G23").Value = 1.246258174420632
$ws.Range("H23").Value = 0.8898665632894165
$ws.Range("L23").Value = 0.5966785001905066
$ws.Range("N23").Value = 1.479682628052899
$ws.Range("B24").Value = 1.472269580428303
$ws.Range("D24").Value = 0.006160236749277459
$ws.Range("E24").Value = 0.6090930465681481
$ws.Range("F24").Value = 1.165226817070334
$ws.Range("G24").Value = 1.113159653110301
$ws.Range("H24").Value = 0.8364645011824621
$ws.Range("L24").Value = 0.5189753720118517
$ws.Range("N24").Value = 1.473536898789988
$ws.Range("B25").Value = 1.343392948847395
$ws.Range("D25").Value = 0.005902821211336118
$ws.Range("E25").Value = 0.5618655337680281
$ws.Range("F25").Value = 1.041501193283167
$ws.Range("G25").Value = 0.9725809287538141
$ws.Range("H25").Value = 0.781235388654494
$ws.Range("L25").Value = 0.4353990564623302
$ws.Range("N25").Value = 1.469393045593975
